$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.02514266666667
$ws.Range("H2").Value = 45.075428
$ws.Range("I2").Value = 0.1401726531301337
$ws.Range("J2").Value = 0.1401726531301337
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 30.99161333333333
$ws.Range("N2").Value = 92.97484
$ws.Range("O2").Value = 0.3599121977633812
$ws.Range("P2").Value = 0.3599121977633811
$ws.Range("Q2").Value = 465.6534118035023
$ws.Range("R2").Value = 4190.88070623152
$ws.Range("S2").Value = 0.05044984765439053
$ws.Range("T2").Value = 0.0504498476543905
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.02514266666667
$ws.Range("H3").Value = 45.075428
$ws.Range("I3").Value = 0.1401726531301337
$ws.Range("J3").Value = 0.1401726531301337
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 29.913269
$ws.Range("N3").Value = 89.739807
$ws.Range("O3").Value = 0.3473891556493311
$ws.Range("P3").Value = 0.3473891556493311
$ws.Range("Q3").Value = 449.4511343513773
$ws.Range("R3").Value = 4045.060209162396
$ws.Range("S3").Value = 0.04869445961600373
$ws.Range("T3").Value = 0.04869445961600371
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.02514266666667
$ws.Range("H4").Value = 45.075428
$ws.Range("I4").Value = 0.1401726531301337
$ws.Range("J4").Value = 0.1401726531301337
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 19.150218
$ws.Range("N4").Value = 57.450654
$ws.Range("O4").Value = 0.2223955550134164
$ws.Range("P4").Value = 0.2223955550134163
$ws.Range("Q4").Value = 287.734757547768
$ws.Range("R4").Value = 2589.612817929912
$ws.Range("S4").Value = 0.03117377499057919
$ws.Range("T4").Value = 0.03117377499057918
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.02514266666667
$ws.Range("H5").Value = 45.075428
$ws.Range("I5").Value = 0.1401726531301337
$ws.Range("J5").Value = 0.1401726531301337
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.053716000000001
$ws.Range("N5").Value = 18.161148
$ws.Range("O5").Value = 0.07030309157387134
$ws.Range("P5").Value = 0.07030309157387132
$ws.Range("Q5").Value = 90.95794656348268
$ws.Range("R5").Value = 818.6215190713441
$ws.Range("S5").Value = 0.009854570869160294
$ws.Range("T5").Value = 0.009854570869160289
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 33.26311566666666
$ws.Range("H6").Value = 99.78934699999999
$ws.Range("I6").Value = 0.3103184627135109
$ws.Range("J6").Value = 0.3103184627135109
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 30.99161333333333
$ws.Range("N6").Value = 92.97484
$ws.Range("O6").Value = 0.3599121977633812
$ws.Range("P6").Value = 0.3599121977633811
$ws.Range("Q6").Value = 1030.877619003275
$ws.Range("R6").Value = 9277.89857102948
$ws.Range("S6").Value = 0.1116873999217736
$ws.Range("T6").Value = 0.1116873999217735
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 33.26311566666666
$ws.Range("H7").Value = 99.78934699999999
$ws.Range("I7").Value = 0.3103184627135109
$ws.Range("J7").Value = 0.3103184627135109
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 29.913269
$ws.Range("N7").Value = 89.739807
$ws.Range("O7").Value = 0.3473891556493311
$ws.Range("P7").Value = 0.3473891556493311
$ws.Range("Q7").Value = 995.0085267151143
$ws.Range("R7").Value = 8955.076740436029
$ws.Range("S7").Value = 0.107801268744445
$ws.Range("T7").Value = 0.107801268744445
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 33.26311566666666
$ws.Range("H8").Value = 99.78934699999999
$ws.Range("I8").Value = 0.3103184627135109
$ws.Range("J8").Value = 0.3103184627135109
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 19.150218
$ws.Range("N8").Value = 57.450654
$ws.Range("O8").Value = 0.2223955550134164
$ws.Range("P8").Value = 0.2223955550134163
$ws.Range("Q8").Value = 636.995916375882
$ws.Range("R8").Value = 5732.963247382938
$ws.Range("S8").Value = 0.06901344674608142
$ws.Range("T8").Value = 0.0690134467460814
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 33.26311566666666
$ws.Range("H9").Value = 99.78934699999999
$ws.Range("I9").Value = 0.3103184627135109
$ws.Range("J9").Value = 0.3103184627135109
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.053716000000001
$ws.Range("N9").Value = 18.161148
$ws.Range("O9").Value = 0.07030309157387134
$ws.Range("P9").Value = 0.07030309157387132
$ws.Range("Q9").Value = 201.3654555211507
$ws.Range("R9").Value = 1812.289099690356
$ws.Range("S9").Value = 0.02181634730121094
$ws.Range("T9").Value = 0.02181634730121093
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 50.12360066666667
$ws.Range("H10").Value = 150.370802
$ws.Range("I10").Value = 0.467613402797773
$ws.Range("J10").Value = 0.4676134027977729
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 30.99161333333333
$ws.Range("N10").Value = 92.97484
$ws.Range("O10").Value = 0.3599121977633812
$ws.Range("P10").Value = 0.3599121977633811
$ws.Range("Q10").Value = 1553.411250735742
$ws.Range("R10").Value = 13980.70125662168
$ws.Range("S10").Value = 0.1682997675045597
$ws.Range("T10").Value = 0.1682997675045596
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 50.12360066666667
$ws.Range("H11").Value = 150.370802
$ws.Range("I11").Value = 0.467613402797773
$ws.Range("J11").Value = 0.4676134027977729
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 29.913269
$ws.Range("N11").Value = 89.739807
$ws.Range("O11").Value = 0.3473891556493311
$ws.Range("P11").Value = 0.3473891556493311
$ws.Range("Q11").Value = 1499.360749990579
$ws.Range("R11").Value = 13494.24674991521
$ws.Range("S11").Value = 0.1624438251682289
$ws.Range("T11").Value = 0.1624438251682289
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 50.12360066666667
$ws.Range("H12").Value = 150.370802
$ws.Range("I12").Value = 0.467613402797773
$ws.Range("J12").Value = 0.4676134027977729
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 19.150218
$ws.Range("N12").Value = 57.450654
$ws.Range("O12").Value = 0.2223955550134164
$ws.Range("P12").Value = 0.2223955550134163
$ws.Range("Q12").Value = 959.877879711612
$ws.Range("R12").Value = 8638.900917404508
$ws.Range("S12").Value = 0.103995142246923
$ws.Range("T12").Value = 0.1039951422469229
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 50.12360066666667
$ws.Range("H13").Value = 150.370802
$ws.Range("I13").Value = 0.467613402797773
$ws.Range("J13").Value = 0.4676134027977729
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.053716000000001
$ws.Range("N13").Value = 18.161148
$ws.Range("O13").Value = 0.07030309157387134
$ws.Range("P13").Value = 0.07030309157387132
$ws.Range("Q13").Value = 303.4340433334107
$ws.Range("R13").Value = 2730.906390000696
$ws.Range("S13").Value = 0.03287466787806142
$ws.Range("T13").Value = 0.0328746678780614
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 8.778397666666669
$ws.Range("H14").Value = 26.335193
$ws.Range("I14").Value = 0.08189548135858246
$ws.Range("J14").Value = 0.08189548135858243
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 30.99161333333333
$ws.Range("N14").Value = 92.97484
$ws.Range("O14").Value = 0.3599121977633812
$ws.Range("P14").Value = 0.3599121977633811
$ws.Range("Q14").Value = 272.056706171569
$ws.Range("R14").Value = 2448.510355544121
$ws.Range("S14").Value = 0.02947518268265743
$ws.Range("T14").Value = 0.02947518268265741
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 8.778397666666669
$ws.Range("H15").Value = 26.335193
$ws.Range("I15").Value = 0.08189548135858246
$ws.Range("J15").Value = 0.08189548135858243
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 29.913269
$ws.Range("N15").Value = 89.739807
$ws.Range("O15").Value = 0.3473891556493311
$ws.Range("P15").Value = 0.3473891556493311
$ws.Range("Q15").Value = 262.5905707919724
$ws.Range("R15").Value = 2363.315137127751
$ws.Range("S15").Value = 0.0284496021206535
$ws.Range("T15").Value = 0.02844960212065348
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 8.778397666666669
$ws.Range("H16").Value = 26.335193
$ws.Range("I16").Value = 0.08189548135858246
$ws.Range("J16").Value = 0.08189548135858243
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 19.150218
$ws.Range("N16").Value = 57.450654
$ws.Range("O16").Value = 0.2223955550134164
$ws.Range("P16").Value = 0.2223955550134163
$ws.Range("Q16").Value = 168.108229007358
$ws.Range("R16").Value = 1512.974061066222
$ws.Range("S16").Value = 0.01821319102983284
$ws.Range("T16").Value = 0.01821319102983283
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 8.778397666666669
$ws.Range("H17").Value = 26.335193
$ws.Range("I17").Value = 0.08189548135858246
$ws.Range("J17").Value = 0.08189548135858243
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.053716000000001
$ws.Range("N17").Value = 18.161148
$ws.Range("O17").Value = 0.07030309157387134
$ws.Range("P17").Value = 0.07030309157387132
$ws.Range("Q17").Value = 53.14192640906268
$ws.Range("R17").Value = 478.2773376815641
$ws.Range("S17").Value = 0.005757505525438696
$ws.Range("T17").Value = 0.005757505525438693
